# RDCC-5182 Added Version check
# Add a new "VERSION" sheet at the end of the workbook containing the
# file version information, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it ends up last.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$versionSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$versionSheet.Name = "VERSION"

# Populate the version info cells.
$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Make B6 the selected cell on the new sheet, and make the new sheet active.
$versionSheet.Range("B6").Select()
